$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 12262.1   # H38: 9433.923000000001 -> 12262.1
$ws.Cells.Item(38, 9).Value = 3880   # I38: 2427.5 -> 3880
$ws.Cells.Item(38, 11).Value = 11640   # K38: 7282.5 -> 11640
$ws.Cells.Item(38, 13).Value = -11268   # M38: -6910.5 -> -11268
$ws.Cells.Item(70, 8).Value = 6166.0835   # H70: 5891.769 -> 6166.0835
$ws.Cells.Item(70, 10).Value = 7742.857   # J70: 7100 -> 7742.857
$ws.Cells.Item(70, 12).Value = 23228.571   # L70: 21300 -> 23228.571
$ws.Cells.Item(70, 14).Value = -23768.571   # N70: -21840 -> -23768.571
$ws.Cells.Item(73, 8).Value = 6166.0835   # H73: 5891.769 -> 6166.0835
$ws.Cells.Item(73, 10).Value = 7742.857   # J73: 7100 -> 7742.857
$ws.Cells.Item(73, 12).Value = 23228.571   # L73: 21300 -> 23228.571
$ws.Cells.Item(73, 14).Value = -25100.571   # N73: -23172 -> -25100.571
$ws.Cells.Item(96, 8).Value = 1380.2858   # H96: 1311.3334 -> 1380.2858
$ws.Cells.Item(96, 9).Value = 630.6923   # I96: 610.3570999999999 -> 630.6923
$ws.Cells.Item(96, 11).Value = 1892.0769   # K96: 1831.0713 -> 1892.0769
$ws.Cells.Item(96, 13).Value = -519.0769   # M96: -458.0712999999998 -> -519.0769
$ws.Cells.Item(98, 8).Value = 1408.1714   # H98: 1498 -> 1408.1714
$ws.Cells.Item(98, 9).Value = 1472.7812   # I98: 1524.4667 -> 1472.7812
$ws.Cells.Item(98, 10).Value = 719   # J98: 704 -> 719
$ws.Cells.Item(98, 11).Value = 1472.7812   # K98: 1524.4667 -> 1472.7812
$ws.Cells.Item(98, 12).Value = 719   # L98: 704 -> 719
$ws.Cells.Item(98, 13).Value = 25.2188000000001   # M98: -26.46669999999995 -> 25.2188000000001
$ws.Cells.Item(98, 14).Value = -3715   # N98: -3700 -> -3715
$ws.Cells.Item(106, 8).Value = 166671680   # H106: 111117784 -> 166671680
$ws.Cells.Item(116, 8).Value = 13547665   # H116: 13217275 -> 13547665
$ws.Cells.Item(116, 9).Value = 16933156   # I116: 16933140 -> 16933156
$ws.Cells.Item(116, 10).Value = 5698.5   # J116: 5309.6665 -> 5698.5
$ws.Cells.Item(116, 11).Value = 16933156   # K116: 16933140 -> 16933156
$ws.Cells.Item(116, 12).Value = 5698.5   # L116: 5309.6665 -> 5698.5
$ws.Cells.Item(116, 13).Value = -16929714   # M116: -16929698 -> -16929714
$ws.Cells.Item(116, 14).Value = -12582.5   # N116: -12193.6665 -> -12582.5
$ws.Cells.Item(122, 8).Value = 1408.1714   # H122: 1498 -> 1408.1714
$ws.Cells.Item(122, 9).Value = 1472.7812   # I122: 1524.4667 -> 1472.7812
$ws.Cells.Item(122, 10).Value = 719   # J122: 704 -> 719
$ws.Cells.Item(122, 11).Value = 4418.3436   # K122: 4573.4001 -> 4418.3436
$ws.Cells.Item(122, 12).Value = 2157   # L122: 2112 -> 2157
$ws.Cells.Item(122, 13).Value = -1968.3436   # M122: -2123.4001 -> -1968.3436
$ws.Cells.Item(122, 14).Value = -7057   # N122: -7012 -> -7057
$ws.Cells.Item(127, 8).Value = 2710.5   # H127: 2722.8518 -> 2710.5
$ws.Cells.Item(127, 9).Value = 1243   # I127: 1273 -> 1243
$ws.Cells.Item(127, 10).Value = 2955.0833   # J127: 2975 -> 2955.0833
$ws.Cells.Item(127, 11).Value = 3729   # K127: 3819 -> 3729
$ws.Cells.Item(127, 12).Value = 8865.249899999999   # L127: 8925 -> 8865.249899999999
$ws.Cells.Item(127, 13).Value = 1231   # M127: 1141 -> 1231
$ws.Cells.Item(127, 14).Value = -18785.2499   # N127: -18845 -> -18785.2499
$ws.Cells.Item(137, 8).Value = 19610298   # H137: 12822641 -> 19610298
$ws.Cells.Item(137, 9).Value = 2312.5   # I137: 1771.8823 -> 2312.5
$ws.Cells.Item(137, 10).Value = 47621704   # J137: 37039836 -> 47621704
$ws.Cells.Item(137, 11).Value = 6937.5   # K137: 5315.6469 -> 6937.5
$ws.Cells.Item(137, 12).Value = 142865112   # L137: 111119508 -> 142865112
$ws.Cells.Item(137, 13).Value = -4387.5   # M137: -2765.6469 -> -4387.5
$ws.Cells.Item(137, 14).Value = -142870212   # N137: -111124608 -> -142870212
$ws.Cells.Item(138, 8).Value = 5028.121   # H138: 5124.5444 -> 5028.121
$ws.Cells.Item(138, 9).Value = 1763.3871   # I138: 1840.1724 -> 1763.3871
$ws.Cells.Item(138, 10).Value = 6714.9   # J138: 6685.9673 -> 6714.9
$ws.Cells.Item(138, 11).Value = 5290.1613   # K138: 5520.5172 -> 5290.1613
$ws.Cells.Item(138, 12).Value = 20144.7   # L138: 20057.9019 -> 20144.7
$ws.Cells.Item(138, 13).Value = -150.1612999999998   # M138: -380.5172000000002 -> -150.1612999999998
$ws.Cells.Item(138, 14).Value = -30424.7   # N138: -30337.9019 -> -30424.7

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 15320.278   # H32: 15815.882 -> 15320.278
$ws.Cells.Item(32, 9).Value = 15095.51   # I32: 15660.039 -> 15095.51
$ws.Cells.Item(32, 11).Value = 15095.51   # K32: 15660.039 -> 15095.51
$ws.Cells.Item(32, 13).Value = -14808.51   # M32: -15373.039 -> -14808.51
$ws.Cells.Item(45, 8).Value = 3384.2354   # H45: 3470.75 -> 3384.2354
$ws.Cells.Item(45, 9).Value = 2925.2727   # I45: 3017.8 -> 2925.2727
$ws.Cells.Item(45, 11).Value = 2925.2727   # K45: 3017.8 -> 2925.2727
$ws.Cells.Item(45, 13).Value = -2548.2727   # M45: -2640.8 -> -2548.2727
$ws.Cells.Item(61, 8).Value = 3795.2144   # H61: 4058.64 -> 3795.2144
$ws.Cells.Item(61, 9).Value = 3359.4348   # I61: 3623.35 -> 3359.4348
$ws.Cells.Item(61, 11).Value = 3359.4348   # K61: 3623.35 -> 3359.4348
$ws.Cells.Item(61, 13).Value = -3147.4348   # M61: -3411.35 -> -3147.4348
$ws.Cells.Item(97, 8).Value = 1618.1818   # H97: 1475 -> 1618.1818
$ws.Cells.Item(97, 9).Value = 787.1429000000001   # I97: 695 -> 787.1429000000001
$ws.Cells.Item(97, 10).Value = 3072.5   # J97: 3035 -> 3072.5
$ws.Cells.Item(97, 11).Value = 787.1429000000001   # K97: 695 -> 787.1429000000001
$ws.Cells.Item(97, 12).Value = 3072.5   # L97: 3035 -> 3072.5
$ws.Cells.Item(97, 13).Value = -291.1429000000001   # M97: -199 -> -291.1429000000001
$ws.Cells.Item(97, 14).Value = -4064.5   # N97: -4027 -> -4064.5
$ws.Cells.Item(122, 8).Value = 4055.9666   # H122: 4151.1035 -> 4055.9666
$ws.Cells.Item(122, 9).Value = 2887.2   # I122: 2953.4583 -> 2887.2
$ws.Cells.Item(122, 11).Value = 8661.599999999999   # K122: 8860.374899999999 -> 8661.599999999999
$ws.Cells.Item(122, 13).Value = -6211.599999999999   # M122: -6410.374899999999 -> -6211.599999999999
$ws.Cells.Item(136, 8).Value = 3795.2144   # H136: 4058.64 -> 3795.2144
$ws.Cells.Item(136, 9).Value = 3359.4348   # I136: 3623.35 -> 3359.4348
$ws.Cells.Item(136, 11).Value = 10078.3044   # K136: 10870.05 -> 10078.3044
$ws.Cells.Item(136, 13).Value = -7528.304400000001   # M136: -8320.049999999999 -> -7528.304400000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1521.48   # H20: 1557.25 -> 1521.48
$ws.Cells.Item(20, 9).Value = 1111.75   # I20: 1152.5454 -> 1111.75
$ws.Cells.Item(20, 11).Value = 1111.75   # K20: 1152.5454 -> 1111.75
$ws.Cells.Item(20, 13).Value = -864.75   # M20: -905.5454 -> -864.75
$ws.Cells.Item(134, 8).Value = 2220.6   # H134: 2672.2856 -> 2220.6
$ws.Cells.Item(134, 9).Value = 2245.111   # I134: 2784.3333 -> 2245.111
$ws.Cells.Item(134, 11).Value = 6735.333   # K134: 8352.999899999999 -> 6735.333
$ws.Cells.Item(134, 13).Value = -4200.333   # M134: -5817.999899999999 -> -4200.333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 34201376   # H132: 35101400 -> 34201376
$ws.Cells.Item(132, 9).Value = 40413148   # I132: 41676044 -> 40413148
$ws.Cells.Item(132, 11).Value = 121239444   # K132: 125028132 -> 121239444
$ws.Cells.Item(132, 13).Value = -121236914   # M132: -125025602 -> -121236914
$ws.Cells.Item(141, 8).Value = 127368.27   # H141: 129245.2 -> 127368.27
$ws.Cells.Item(141, 10).Value = 127368.27   # J141: 129245.2 -> 127368.27
$ws.Cells.Item(141, 12).Value = 127368.27   # L141: 129245.2 -> 127368.27
$ws.Cells.Item(141, 14).Value = -137728.27   # N141: -139605.2 -> -137728.27

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(70, 8).Value = 9796.625   # H70: 333350000 -> 9796.625
$ws.Cells.Item(70, 10).Value = 9796.625   # J70: 333350000 -> 9796.625
$ws.Cells.Item(70, 12).Value = 29389.875   # L70: 1000050000 -> 29389.875
$ws.Cells.Item(70, 14).Value = -30019.875   # N70: -1000050630 -> -30019.875
$ws.Cells.Item(73, 8).Value = 9796.625   # H73: 333350000 -> 9796.625
$ws.Cells.Item(73, 10).Value = 9796.625   # J73: 333350000 -> 9796.625
$ws.Cells.Item(73, 12).Value = 29389.875   # L73: 1000050000 -> 29389.875
$ws.Cells.Item(73, 14).Value = -31573.875   # N73: -1000052184 -> -31573.875
$ws.Cells.Item(131, 8).Value = 15937474   # H131: 16506605 -> 15937474
$ws.Cells.Item(131, 10).Value = 16471678   # J131: 17220308 -> 16471678
$ws.Cells.Item(131, 12).Value = 49415034   # L131: 51660924 -> 49415034
$ws.Cells.Item(131, 14).Value = -49425114   # N131: -51671004 -> -49425114

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 1864055.1   # H80: 1863777.2 -> 1864055.1
$ws.Cells.Item(80, 9).Value = 3329299.5   # I80: 2774832.8 -> 3329299.5
$ws.Cells.Item(80, 10).Value = 32499.75   # J80: 41666.332 -> 32499.75
$ws.Cells.Item(80, 11).Value = 3329299.5   # K80: 2774832.8 -> 3329299.5
$ws.Cells.Item(80, 12).Value = 32499.75   # L80: 41666.332 -> 32499.75
$ws.Cells.Item(80, 13).Value = -3328301.5   # M80: -2773834.8 -> -3328301.5
$ws.Cells.Item(80, 14).Value = -34495.75   # N80: -43662.332 -> -34495.75
$ws.Cells.Item(83, 8).Value = 1864055.1   # H83: 1863777.2 -> 1864055.1
$ws.Cells.Item(83, 9).Value = 3329299.5   # I83: 2774832.8 -> 3329299.5
$ws.Cells.Item(83, 10).Value = 32499.75   # J83: 41666.332 -> 32499.75
$ws.Cells.Item(83, 11).Value = 16646497.5   # K83: 13874164 -> 16646497.5
$ws.Cells.Item(83, 12).Value = 162498.75   # L83: 208331.66 -> 162498.75
$ws.Cells.Item(83, 13).Value = -16641505.5   # M83: -13869172 -> -16641505.5
$ws.Cells.Item(83, 14).Value = -172482.75   # N83: -218315.66 -> -172482.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1631   # H61: 1621.6 -> 1631
$ws.Cells.Item(61, 9).Value = 1631   # I61: 1621.6 -> 1631
$ws.Cells.Item(61, 11).Value = 1631   # K61: 1621.6 -> 1631
$ws.Cells.Item(61, 13).Value = -1429   # M61: -1419.6 -> -1429
$ws.Cells.Item(113, 8).Value = 1631   # H113: 1621.6 -> 1631
$ws.Cells.Item(113, 9).Value = 1631   # I113: 1621.6 -> 1631
$ws.Cells.Item(113, 11).Value = 1631   # K113: 1621.6 -> 1631
$ws.Cells.Item(113, 13).Value = 539   # M113: 548.4000000000001 -> 539
$ws.Cells.Item(136, 8).Value = 4539   # H136: 4826.9443 -> 4539
$ws.Cells.Item(136, 9).Value = 2308.5557   # I136: 2411.7144 -> 2308.5557
$ws.Cells.Item(136, 11).Value = 6925.6671   # K136: 7235.1432 -> 6925.6671
$ws.Cells.Item(136, 13).Value = -4375.6671   # M136: -4685.1432 -> -4375.6671
$ws.Cells.Item(139, 8).Value = 100715   # H139: 0 -> 100715
$ws.Cells.Item(139, 10).Value = 100715   # J139: 0 -> 100715
$ws.Cells.Item(139, 12).Value = 100715   # L139: 0 -> 100715
$ws.Cells.Item(139, 14).Value = -110995   # N139: None -> -110995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3477137.8   # H81: 3477037.8 -> 3477137.8
$ws.Cells.Item(81, 9).Value = 10417416   # I81: 6945410.5 -> 10417416
$ws.Cells.Item(81, 10).Value = 6998.75   # J81: 8665 -> 6998.75
$ws.Cells.Item(81, 11).Value = 20834832   # K81: 13890821 -> 20834832
$ws.Cells.Item(81, 12).Value = 13997.5   # L81: 17330 -> 13997.5
$ws.Cells.Item(81, 13).Value = -20833771   # M81: -13889760 -> -20833771
$ws.Cells.Item(81, 14).Value = -16119.5   # N81: -19452 -> -16119.5
$ws.Cells.Item(84, 8).Value = 3477137.8   # H84: 3477037.8 -> 3477137.8
$ws.Cells.Item(84, 9).Value = 10417416   # I84: 6945410.5 -> 10417416
$ws.Cells.Item(84, 10).Value = 6998.75   # J84: 8665 -> 6998.75
$ws.Cells.Item(84, 11).Value = 104174160   # K84: 69454105 -> 104174160
$ws.Cells.Item(84, 12).Value = 69987.5   # L84: 86650 -> 69987.5
$ws.Cells.Item(84, 13).Value = -104168856   # M84: -69448801 -> -104168856
$ws.Cells.Item(84, 14).Value = -80595.5   # N84: -97258 -> -80595.5
$ws.Cells.Item(126, 8).Value = 6205.375   # H126: 6353.4 -> 6205.375
$ws.Cells.Item(126, 9).Value = 5134   # I126: 5298.143 -> 5134
$ws.Cells.Item(126, 11).Value = 15402   # K126: 15894.429 -> 15402
$ws.Cells.Item(126, 13).Value = -12932   # M126: -13424.429 -> -12932
$ws.Cells.Item(132, 8).Value = 55560924   # H132: 62505764 -> 55560924
$ws.Cells.Item(132, 9).Value = 3051.5   # I132: 3904 -> 3051.5
$ws.Cells.Item(132, 11).Value = 9154.5   # K132: 11712 -> 9154.5
$ws.Cells.Item(132, 13).Value = -6624.5   # M132: -9182 -> -6624.5
